$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 29 (shifts old row 29 -> row 30)
$ws.Rows("29:29").Insert()

# Fill in the new row 29 with the updated weekly entry
$ws.Range("A29").Value = 10
$ws.Range("B29").Value = "Vega Modelo de Temuco"
$ws.Range("C29").Value = "La Araucanía"
$ws.Range("D29").Value = 45075
$ws.Range("D29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E29").Value = 9
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100108
$ws.Range("H29").Value = "Tropicales y subtropicales"
$ws.Range("I29").Value = 100108001
$ws.Range("J29").Value = "Guayaba"
$ws.Range("K29").Value = "Sin especificar"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 240
$ws.Range("N29").Value = 3200
$ws.Range("O29").Value = 3200
$ws.Range("P29").Value = 3200
$ws.Range("Q29").Value = "`$/kilo"
$ws.Range("R29").Value = "Región de Arica y Parinacota"
$ws.Range("S29").Value = 3200
$ws.Range("T29").Value = 1
